$p = $ppt.ActivePresentation
$s = $p.Slides.Item(35)

# --- Content Placeholder 2 (left column): swap "Blog" text for "Python notebook" text ---
$shape3 = $s.Shapes.Item(2)
$tr3 = $shape3.TextFrame.TextRange
$tr3.Paragraphs(3).Runs(1).Text = "Python notebook: excellent"
$tr3.Paragraphs(4).Runs(1).Text = "Quarto adds to Python notebooks without detracting anything. All you need are a few YAML lines."

# --- Content Placeholder 3 (right column): swap "Python notebook" text for "Blog" text ---
$shape4 = $s.Shapes.Item(3)
$tr4 = $shape4.TextFrame.TextRange
$tr4.Paragraphs(1).Runs(1).Text = "Blog: excellent"
$tr4.Paragraphs(2).Runs(1).Text = "Quarto allows me to have a scriptable, Python-based blog. I can automate my blog to tweet and post to LinkedIn when I write new articles."

# --- Split "Presentations: only if you have" into its own header + a plain follow-up line ---
# Insert two brand-new paragraphs (cloning the plain "buNone" formatting used by the
# paragraph right above) straight after the "Quarto adds..." paragraph, then move the
# bold header formatting onto the first of the two, and drop the paragraph it replaces.
$para2 = $tr4.Paragraphs(2)
$para2.InsertAfter([char]13 + "Presentations: unsure" + [char]13 + "Only if you have") | Out-Null

$tr4b = $shape4.TextFrame.TextRange
$newHeader = $tr4b.Paragraphs(3)
$newHeader.ParagraphFormat.SpaceBefore = 30
$newHeader.ParagraphFormat.Bullet.Visible = -1
$newHeader.ParagraphFormat.Bullet.Visible = 0
$newHeader.Font.Bold = -1

$tr4c = $shape4.TextFrame.TextRange
$tr4c.Paragraphs(5).Delete()
